$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 7).Value = 7.871715666666666
$ws.Cells.Item(2, 8).Value = 23.615147
$ws.Cells.Item(2, 9).Value = 0.02771913691218268
$ws.Cells.Item(2, 10).Value = 0.02771913691218268
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 13).Value = 7.214110666666667
$ws.Cells.Item(2, 14).Value = 21.642332
$ws.Cells.Item(2, 15).Value = 0.4688823795981188
$ws.Cells.Item(2, 16).Value = 0.4688823795981188
$ws.Cells.Item(2, 17).Value = 56.7874279558671
$ws.Cells.Item(2, 18).Value = 511.0868516028039
$ws.Cells.Item(2, 19).Value = 0.01299701487579027
$ws.Cells.Item(2, 20).Value = 0.01299701487579027

$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 7).Value = 7.871715666666666
$ws.Cells.Item(3, 8).Value = 23.615147
$ws.Cells.Item(3, 9).Value = 0.02771913691218268
$ws.Cells.Item(3, 10).Value = 0.02771913691218268
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 13).Value = 7.110350666666666
$ws.Cells.Item(3, 14).Value = 21.331052
$ws.Cells.Item(3, 15).Value = 0.4621384803214003
$ws.Cells.Item(3, 16).Value = 0.4621384803214003
$ws.Cells.Item(3, 17).Value = 55.97065873829377
$ws.Cells.Item(3, 18).Value = 503.7359286446439
$ws.Cells.Item(3, 19).Value = 0.01281007980841694
$ws.Cells.Item(3, 20).Value = 0.01281007980841694

$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 7).Value = 7.871715666666666
$ws.Cells.Item(4, 8).Value = 23.615147
$ws.Cells.Item(4, 9).Value = 0.02771913691218268
$ws.Cells.Item(4, 10).Value = 0.02771913691218268
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 13).Value = 1.061296333333333
$ws.Cells.Item(4, 14).Value = 3.183889
$ws.Cells.Item(4, 15).Value = 0.06897914008048092
$ws.Cells.Item(4, 16).Value = 0.06897914008048092
$ws.Cells.Item(4, 17).Value = 8.354222974075887
$ws.Cells.Item(4, 18).Value = 75.18800676668299
$ws.Cells.Item(4, 19).Value = 0.001912042227975479
$ws.Cells.Item(4, 20).Value = 0.001912042227975478

$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 7).Value = 266.1315866666666
$ws.Cells.Item(5, 8).Value = 798.3947599999999
$ws.Cells.Item(5, 9).Value = 0.9371448614065047
$ws.Cells.Item(5, 10).Value = 0.9371448614065045
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 13).Value = 7.214110666666667
$ws.Cells.Item(5, 14).Value = 21.642332
$ws.Cells.Item(5, 15).Value = 0.4688823795981188
$ws.Cells.Item(5, 16).Value = 0.4688823795981188
$ws.Cells.Item(5, 17).Value = 1919.902718108924
$ws.Cells.Item(5, 18).Value = 17279.12446298032
$ws.Cells.Item(5, 19).Value = 0.4394107126444312
$ws.Cells.Item(5, 20).Value = 0.4394107126444312

$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 7).Value = 266.1315866666666
$ws.Cells.Item(6, 8).Value = 798.3947599999999
$ws.Cells.Item(6, 9).Value = 0.9371448614065047
$ws.Cells.Item(6, 10).Value = 0.9371448614065045
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 13).Value = 7.110350666666666
$ws.Cells.Item(6, 14).Value = 21.331052
$ws.Cells.Item(6, 15).Value = 0.4621384803214003
$ws.Cells.Item(6, 16).Value = 0.4621384803214003
$ws.Cells.Item(6, 17).Value = 1892.288904676391
$ws.Cells.Item(6, 18).Value = 17030.60014208752
$ws.Cells.Item(6, 19).Value = 0.4330907020914114
$ws.Cells.Item(6, 20).Value = 0.4330907020914113

$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 7).Value = 266.1315866666666
$ws.Cells.Item(7, 8).Value = 798.3947599999999
$ws.Cells.Item(7, 9).Value = 0.9371448614065047
$ws.Cells.Item(7, 10).Value = 0.9371448614065045
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 13).Value = 1.061296333333333
$ws.Cells.Item(7, 14).Value = 3.183889
$ws.Cells.Item(7, 15).Value = 0.06897914008048092
$ws.Cells.Item(7, 16).Value = 0.06897914008048092
$ws.Cells.Item(7, 17).Value = 282.4444771135155
$ws.Cells.Item(7, 18).Value = 2542.00029402164
$ws.Cells.Item(7, 19).Value = 0.06464344667066216
$ws.Cells.Item(7, 20).Value = 0.06464344667066216

$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 7).Value = 9.977966333333333
$ws.Cells.Item(8, 8).Value = 29.933899
$ws.Cells.Item(8, 9).Value = 0.03513600168131278
$ws.Cells.Item(8, 10).Value = 0.03513600168131277
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 13).Value = 7.214110666666667
$ws.Cells.Item(8, 14).Value = 21.642332
$ws.Cells.Item(8, 15).Value = 0.4688823795981188
$ws.Cells.Item(8, 16).Value = 0.4688823795981188
$ws.Cells.Item(8, 17).Value = 71.98215335694088
$ws.Cells.Item(8, 18).Value = 647.839380212468
$ws.Cells.Item(8, 19).Value = 0.01647465207789744
$ws.Cells.Item(8, 20).Value = 0.01647465207789744

$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 7).Value = 9.977966333333333
$ws.Cells.Item(9, 8).Value = 29.933899
$ws.Cells.Item(9, 9).Value = 0.03513600168131278
$ws.Cells.Item(9, 10).Value = 0.03513600168131277
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 13).Value = 7.110350666666666
$ws.Cells.Item(9, 14).Value = 21.331052
$ws.Cells.Item(9, 15).Value = 0.4621384803214003
$ws.Cells.Item(9, 16).Value = 0.4621384803214003
$ws.Cells.Item(9, 17).Value = 70.94683957019421
$ws.Cells.Item(9, 18).Value = 638.521556131748
$ws.Cells.Item(9, 19).Value = 0.01623769842157205
$ws.Cells.Item(9, 20).Value = 0.01623769842157205

$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 7).Value = 9.977966333333333
$ws.Cells.Item(10, 8).Value = 29.933899
$ws.Cells.Item(10, 9).Value = 0.03513600168131278
$ws.Cells.Item(10, 10).Value = 0.03513600168131277
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 13).Value = 1.061296333333333
$ws.Cells.Item(10, 14).Value = 3.183889
$ws.Cells.Item(10, 15).Value = 0.06897914008048092
$ws.Cells.Item(10, 16).Value = 0.06897914008048092
$ws.Cells.Item(10, 17).Value = 10.58957908369011
$ws.Cells.Item(10, 18).Value = 95.30621175321099
$ws.Cells.Item(10, 19).Value = 0.002423651181843287
$ws.Cells.Item(10, 20).Value = 0.002423651181843287
